$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'242.66"
$ws.Range("G2").Value = "'19"
$ws.Range("D3").Value = "'23.53"
$ws.Range("G3").Value = "'19"
$ws.Range("D4").Value = "'5.657"
$ws.Range("G4").Value = "'19"
$ws.Range("D5").Value = "'0.05813"
$ws.Range("G5").Value = "'19"
$ws.Range("D6").Value = "'3.411"
$ws.Range("G6").Value = "'19"
$ws.Range("D7").Value = "'6.469"
$ws.Range("G7").Value = "'19"
$ws.Range("D8").Value = "'1.318"
$ws.Range("G8").Value = "'19"
$ws.Range("D9").Value = "'0.7974"
$ws.Range("G9").Value = "'19"
$ws.Range("D10").Value = "'0.1459"
$ws.Range("G10").Value = "'19"
$ws.Range("D11").Value = "'0.07626"
$ws.Range("G11").Value = "'19"
$ws.Range("D12").Value = "'0.03256"
$ws.Range("G12").Value = "'19"
$ws.Range("D13").Value = "'0.02939"
$ws.Range("G13").Value = "'19"
$ws.Range("D14").Value = "'0.09227"
$ws.Range("G14").Value = "'19"
$ws.Range("D15").Value = "'0.001659"
$ws.Range("G15").Value = "'19"
$ws.Range("G16").Value = "'19"
$ws.Range("D17").Value = "'0.04745"
$ws.Range("G17").Value = "'19"
$ws.Range("D18").Value = "'0.0005992"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("G18").Value = "'19"
$ws.Range("G19").Value = "'19"
$ws.Range("D20").Value = "'0.005465"
$ws.Range("G20").Value = "'19"
$ws.Range("G21").Value = "'19"
$ws.Range("D22").Value = "'0.0001501"
$ws.Range("G22").Value = "'19"
$ws.Range("D23").Value = "'3.696"
$ws.Range("G23").Value = "'19"
$ws.Range("D24").Value = "'2.175"
$ws.Range("G24").Value = "'19"
$ws.Range("D25").Value = "'0.3328"
$ws.Range("G25").Value = "'19"
$ws.Range("D26").Value = "'0.1235"
$ws.Range("G26").Value = "'19"
$ws.Range("D27").Value = "'0.001000"
$ws.Range("E27").Value = "26UpBotsUBXT"
$ws.Range("G27").Value = "'19"
$ws.Range("G28").Value = "'19"
$ws.Range("G29").Value = "'19"
$ws.Range("G30").Value = "'19"
$ws.Range("G31").Value = "'19"
$ws.Range("G32").Value = "'19"
$ws.Range("G33").Value = "'19"
$ws.Range("G34").Value = "'19"
$ws.Range("G35").Value = "'19"
$ws.Range("G36").Value = "'19"
$ws.Range("G37").Value = "'19"
$ws.Range("G38").Value = "'19"
$ws.Range("G39").Value = "'19"
$ws.Range("D40").Value = "'0.04298"
$ws.Range("G40").Value = "'19"
$ws.Range("D41").Value = "'0.007122"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("G41").Value = "'19"
$ws.Range("G42").Value = "'19"
$ws.Range("D43").Value = "'0.1053"
$ws.Range("G43").Value = "'19"
$ws.Range("D44").Value = "'0.009534"
$ws.Range("G44").Value = "'19"
$ws.Range("E45").Value = "44ACDXExchangeACXTWorstin24h"
$ws.Range("G45").Value = "'19"
$ws.Range("D46").Value = "'0.00005361"
$ws.Range("G46").Value = "'19"
$ws.Range("G47").Value = "'19"
$ws.Range("D48").Value = "'0.7856"
$ws.Range("G48").Value = "'19"
$ws.Range("D49").Value = "'0.1027"
$ws.Range("G49").Value = "'19"
$ws.Range("G50").Value = "'19"
$ws.Range("G51").Value = "'19"
